# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Burundi" and "San Cristobal y Nieves" rows (18/19 listing order) ---
# Row 198 currently holds "San Cristobal y Nieves" data; Row 199 holds "Burundi" data.
# After the edit, row 198 should show Burundi (with Burundi's stats) and row 199
# should show San Cristobal y Nieves (with its stats) - i.e. the two countries swap places.
$ws.Range("A198").Value = "Burundi"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 7
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 8
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# --- Row 18: Paises Bajos (Netherlands) stats update ---
$ws.Range("B18").Value = 40571
$ws.Range("C18").Value = 335
$ws.Range("E18").Value = 35265
$ws.Range("F18").Value = 861
$ws.Range("G18").Value = 69
$ws.Range("H18").Value = 5056

# --- Row 60: Kazajistan stats update ---
$ws.Range("B60").Value = 3877
$ws.Range("C60").Value = 20
$ws.Range("E60").Value = 2857
$ws.Range("F60").Value = 40

# --- Row 110: Georgia stats update ---
$ws.Range("D110").Value = 221
$ws.Range("E110").Value = 359
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 9
